$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values (e.g. "335.08", "10.10")
# as plain text in the original workbook, preserving exact digits/trailing
# zeros. A bare .Value assignment of such a string gets auto-coerced to a
# floating point number by Excel, so for any replacement price that parses
# as a number we briefly force the cell to Text format, assign the literal
# string, then restore the default style (no explicit style index), which
# is how the original cells are stored.

$ws.Range("D2").Value = '28.565.39'
$ws.Range("E2").Value = '  -3.57%  '
$ws.Range("D3").Value = '1.847.95'
$ws.Range("E3").Value = '  -4.14%  '
$ws.Range("E4").Value = '  -1.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4644'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3898'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07896'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9754'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.33%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.807'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.814.80'
$ws.Range("E14").Value = '  -3.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.959'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06902'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.76'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.003'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Value = '28.587.92'
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.382'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.152'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.84%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.35%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.038'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.05%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.004'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.94%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9651'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.63%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09358'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.363'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.70%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.463'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.345'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.25%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06100'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.78%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02197'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.55%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.165'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5694'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.97%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.656'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.51%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.89%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1791'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.422'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.250'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5372'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07090'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.05%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.899'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.58%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.345'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.62%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.14%  '
